$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. Create the two new market sheets ("Norway", "Poland") by copying the
#    "Greece" market sheet, which carries the exact formatting (styles,
#    column widths, merged cells, page setup) the new sheets should inherit.
#    The copies are placed at the end of the tab order (after "Hungary").
# ---------------------------------------------------------------------------
$greece = $wb.Worksheets.Item("Greece")
$hungary = $wb.Worksheets.Item("Hungary")

$greece.Copy($null, $hungary)
$norway = $wb.Worksheets.Item($wb.Worksheets.Count)
$norway.Name = "Norway"
$norway.Range("B4").Value = "NGC-2931/T3061"
$norway.Range("B2").Value = "Norway Market"

$greece.Copy($null, $norway)
$poland = $wb.Worksheets.Item($wb.Worksheets.Count)
$poland.Name = "Poland"
$poland.Range("B4").Value = "NGC-2920/T3104"
$poland.Range("B2").Value = "Poland Market"

# ---------------------------------------------------------------------------
# 2. Add the two new accessory rows ("MX-DPBX" / "MX-BBX") to every market
#    sheet's accessory list (rows 9-10, pushing the previous rows 9-12 down
#    to rows 11-14). Insert format is copied from the existing row 8 so the
#    borders match the rest of the list.
# ---------------------------------------------------------------------------
function Add-AccessoryRows {
    param($ws, $firstValue, $secondValue)

    $ws.Rows("9:10").Insert()
    $ws.Range("A8").Copy()
    $ws.Range("A9:A10").PasteSpecial(-4122)
    $ws.Range("A9").Value = $firstValue
    $ws.Range("A10").Value = $secondValue
}

$portugal = $wb.Worksheets.Item("Portugal")
Add-AccessoryRows $portugal "MX-DPBX" "MX-BBX"

$croatia = $wb.Worksheets.Item("Croatia")
Add-AccessoryRows $croatia "MX-BBX" "MX-DPBX"

Add-AccessoryRows $greece "MX-DPBX" "MX-BBX"

Add-AccessoryRows $norway "MX-BBX" "MX-DPBX"
Add-AccessoryRows $poland "MX-DPBX" "MX-BBX"

# ---------------------------------------------------------------------------
# 3. Update the selection / active-cell state on each touched sheet, ending
#    with "Poland" as the active tab (matching the new activeTab).
# ---------------------------------------------------------------------------
$portugal.Activate()
$portugal.Range("A10").Select()

$croatia.Activate()
$croatia.Range("A10").Select()

$greece.Activate()
$greece.Range("A10").Select()

$norway.Activate()
$norway.Range("A10").Select()

$poland.Activate()
$poland.Range("A10").Select()
